$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix "Webdev" casing -> "WebDev" for all existing rows that have it ---
$ws.Range("B2").Value = "WebDev"
$ws.Range("B4").Value = "WebDev"
$ws.Range("B6").Value = "WebDev"
$ws.Range("B7").Value = "WebDev"
$ws.Range("B10").Value = "WebDev"
$ws.Range("B12").Value = "WebDev"
$ws.Range("B14").Value = "WebDev"
$ws.Range("B18").Value = "WebDev"
$ws.Range("B20").Value = "WebDev"
$ws.Range("B22").Value = "WebDev"
$ws.Range("B23").Value = "WebDev"
$ws.Range("B30").Value = "WebDev"
$ws.Range("B31").Value = "WebDev"

# --- Step 2: new people, rows 32-40 (Nome31..Nome39), column A first ---
$ws.Range("A32").Value = "Nome31"
$ws.Range("A33").Value = "Nome32"
$ws.Range("A34").Value = "Nome33"
$ws.Range("A35").Value = "Nome34"
$ws.Range("A36").Value = "Nome35"
$ws.Range("A37").Value = "Nome36"
$ws.Range("A38").Value = "Nome37"
$ws.Range("A39").Value = "Nome38"
$ws.Range("A40").Value = "Nome39"

# --- Step 3: teams for rows 32-40 (Business for row 33 deferred to the end) ---
$ws.Range("B32").Value = "WebDev"
$ws.Range("B36").Value = "Logistics"
$ws.Range("B37").Value = "Logistics"
$ws.Range("B38").Value = "Logistics"
$ws.Range("B39").Value = "Logistics"
$ws.Range("B34").Value = "Coordinator"
$ws.Range("B35").Value = "Coordinator"
$ws.Range("B40").Value = "Speakers"

# --- Step 4: new people, rows 41-52 (Nome40..Nome51), column A ---
$ws.Range("A41").Value = "Nome40"
$ws.Range("A42").Value = "Nome41"
$ws.Range("A43").Value = "Nome42"
$ws.Range("A44").Value = "Nome43"
$ws.Range("A45").Value = "Nome44"
$ws.Range("A46").Value = "Nome45"
$ws.Range("A47").Value = "Nome46"
$ws.Range("A48").Value = "Nome47"
$ws.Range("A49").Value = "Nome48"
$ws.Range("A50").Value = "Nome49"
$ws.Range("A51").Value = "Nome50"
$ws.Range("A52").Value = "Nome51"

# --- Step 5: teams for rows 41-52 (Business entries deferred to the end) ---
$ws.Range("B41").Value = "Coordinator"
$ws.Range("B42").Value = "Logistics"
$ws.Range("B43").Value = "Logistics"
$ws.Range("B45").Value = "Coordinator"
$ws.Range("B46").Value = "Logistics"
$ws.Range("B48").Value = "Logistics"
$ws.Range("B51").Value = "Logistics"
$ws.Range("B52").Value = "Logistics"

# --- Step 6: column C (availability) for the new rows ---
$avail = "[1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1]"
$ws.Range("C32").Value = $avail
$ws.Range("C33").Value = $avail
$ws.Range("C34").Value = $avail
$ws.Range("C35").Value = $avail
$ws.Range("C36").Value = $avail
$ws.Range("C37").Value = $avail
$ws.Range("C38").Value = $avail
$ws.Range("C39").Value = $avail
$ws.Range("C40").Value = $avail
$ws.Range("C41").Value = $avail
$ws.Range("C42").Value = $avail
$ws.Range("C43").Value = $avail
$ws.Range("C44").Value = $avail
$ws.Range("C45").Value = $avail
$ws.Range("C46").Value = $avail
$ws.Range("C47").Value = $avail
$ws.Range("C48").Value = $avail
$ws.Range("C49").Value = $avail
$ws.Range("C50").Value = $avail
$ws.Range("C51").Value = $avail
$ws.Range("C52").Value = $avail

# --- Step 7: every "Business" assignment last, so the new shared string sorts after everything else ---
$ws.Range("B24").Value = "Business"
$ws.Range("B25").Value = "Business"
$ws.Range("B28").Value = "Business"
$ws.Range("B33").Value = "Business"
$ws.Range("B44").Value = "Business"
$ws.Range("B47").Value = "Business"
$ws.Range("B49").Value = "Business"
$ws.Range("B50").Value = "Business"

# --- Step 8: view state - scroll so the bottom of the new data is visible, selection on B45 ---
$ws.Range("B45").Select()
$excel.ActiveWindow.ScrollRow = 29

Write-Output "edit applied"
